$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9714423418045044
$ws.Range("B1").Value = 1.028986811637878
$ws.Range("C1").Value = 5.117161750793457
$ws.Range("D1").Value = 2.068168878555298
$ws.Range("E1").Value = 1.215786457061768
